$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two example address rows with real property names ---
$ws.Range("A2").Value = "Panattoni Park Cheb"
$ws.Range("A3").Value = "GLP Budapest"

# --- Explicitly (re-)lock the header row cells (Format Cells > Protection > Locked) ---
$ws.Range("A1:V1").Locked = $true

# --- Update page setup (paper size / orientation) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Move the active selection ---
$ws.Range("I21").Select()
